# Update "想去人数" (number of people interested) counts in column F
# for rows 2,3,5,8,9,10 on both the "展览" and "全部类型" worksheets.
# These two sheets contain duplicate data sets, and both need the same update.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 614
    3  = 434
    5  = 19
    8  = 1076
    9  = 3867
    10 = 77
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
